$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$afterTitle = $titlePara.Range.Duplicate
$afterTitle.Collapse(0)               # wdCollapseEnd
$afterTitle.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"      # plain paragraph, no heading style

$metaBold = "Meta description"
$metaRest = ": Experience the thrill of Congo King Quad Shot - a highly volatile slot game with captivating graphics and chance to win one of four jackpots. Play for free now!"
$metaPara.Range.Text = $metaBold + $metaRest

# Make "Meta description" bold.
$boldRange = $metaPara.Range.Duplicate
$boldRange.Start = $metaPara.Range.Start
$boldRange.End = $metaPara.Range.Start + $metaBold.Length
$boldRange.Bold = 1

# Leading empty run, matching the pattern used by the rest of the doc.
$leadIn = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)
$leadIn.InsertBefore("")

# ------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph that used to sit near
#    the end of the document (just before the italic meta paragraph).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the old meta-description sentence (now the very last
#    paragraph, italic) with the new DALL-E image prompt text.
# ------------------------------------------------------------------
$oldText = "Experience the thrill of Congo King Quad Shot - a highly volatile slot game with captivating graphics and chance to win one of four jackpots. Play for free now!"
$newText = 'DALLE, please create a feature image fitting the game "Congo King Quad Shot" that meets the following requirements: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses. The image should capture the adventurous spirit of the game and convey the excitement of exploring through the jungle. Please ensure the colors used in the image are vibrant and eye-catching. The image should be appealing and encourage players to take a chance on the game.'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

Write-Output "done"
